# Applies numeric-value updates to H:N columns across several sheets,
# matching a scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2063.9285
$ws.Range("I51").Value = 1868.75
$ws.Range("J51").Value = 2324.1667
$ws.Range("K51").Value = 1868.75
$ws.Range("L51").Value = 2324.1667
$ws.Range("M51").Value = -1384.75
$ws.Range("N51").Value = -3292.1667

$ws.Range("H132").Value = 22224778
$ws.Range("I132").Value = 28573002
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 85719006
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -85716476
$ws.Range("N132").Value = -23060

$ws.Range("H135").Value = 1481.6364
$ws.Range("I135").Value = 1785.5714
$ws.Range("J135").Value = 949.75
$ws.Range("K135").Value = 16070.1426
$ws.Range("L135").Value = 8547.75
$ws.Range("M135").Value = -13535.1426
$ws.Range("N135").Value = -13617.75

$ws.Range("H137").Value = 5267966
$ws.Range("I137").Value = 10006700
$ws.Range("J137").Value = 2705.889
$ws.Range("K137").Value = 30020100
$ws.Range("L137").Value = 8117.667
$ws.Range("M137").Value = -30017550
$ws.Range("N137").Value = -13217.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 9261306
$ws.Range("I2").Value = 13159014
$ws.Range("J2").Value = 4249.75
$ws.Range("K2").Value = 13159014
$ws.Range("L2").Value = 4249.75
$ws.Range("M2").Value = -13158901
$ws.Range("N2").Value = -4475.75

$ws.Range("H45").Value = 1722.125
$ws.Range("I45").Value = 1099.2
$ws.Range("J45").Value = 3946.8572
$ws.Range("K45").Value = 1099.2
$ws.Range("L45").Value = 3946.8572
$ws.Range("M45").Value = -722.2
$ws.Range("N45").Value = -4700.8572

$ws.Range("H110").Value = 1239.619
$ws.Range("I110").Value = 688.7143
$ws.Range("K110").Value = 688.7143
$ws.Range("M110").Value = 1356.2857

$ws.Range("H116").Value = 9261306
$ws.Range("I116").Value = 13159014
$ws.Range("J116").Value = 4249.75
$ws.Range("K116").Value = 13159014
$ws.Range("L116").Value = 4249.75
$ws.Range("M116").Value = -13156720
$ws.Range("N116").Value = -8837.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9261306
$ws.Range("I3").Value = 13159014
$ws.Range("J3").Value = 4249.75
$ws.Range("K3").Value = 13159014
$ws.Range("L3").Value = 4249.75
$ws.Range("M3").Value = -13158900
$ws.Range("N3").Value = -4477.75

$ws.Range("H36").Value = 35539
$ws.Range("I36").Value = 1037
$ws.Range("J36").Value = 70041
$ws.Range("K36").Value = 1037
$ws.Range("L36").Value = 70041
$ws.Range("M36").Value = -503
$ws.Range("N36").Value = -71109

$ws.Range("H86").Value = 2091.0435
$ws.Range("I86").Value = 1466.3334
$ws.Range("J86").Value = 4340
$ws.Range("K86").Value = 1466.3334
$ws.Range("L86").Value = 4340
$ws.Range("M86").Value = -343.3334
$ws.Range("N86").Value = -6586

$ws.Range("H89").Value = 2091.0435
$ws.Range("I89").Value = 1466.3334
$ws.Range("J89").Value = 4340
$ws.Range("K89").Value = 7331.666999999999
$ws.Range("L89").Value = 21700
$ws.Range("M89").Value = -1715.666999999999
$ws.Range("N89").Value = -32932

$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws.Range("H133").Value = 40000
$ws.Range("J133").Value = 40000
$ws.Range("L133").Value = 40000
$ws.Range("N133").Value = -50120

$ws.Range("H134").Value = 1962.8334
$ws.Range("I134").Value = 1962.8334
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5888.5002
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -3353.5002
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3228180.5
$ws.Range("I31").Value = 3573317.5
$ws.Range("J31").Value = 6900
$ws.Range("K31").Value = 3573317.5
$ws.Range("L31").Value = 6900
$ws.Range("M31").Value = -3573022.5
$ws.Range("N31").Value = -7490

$ws.Range("H34").Value = 3228180.5
$ws.Range("I34").Value = 3573317.5
$ws.Range("J34").Value = 6900
$ws.Range("K34").Value = 3573317.5
$ws.Range("L34").Value = 6900
$ws.Range("M34").Value = -3573115.5
$ws.Range("N34").Value = -7304

$ws.Range("H99").Value = 2152.3125
$ws.Range("I99").Value = 2148.3333
$ws.Range("J99").Value = 2154.7
$ws.Range("K99").Value = 2148.3333
$ws.Range("L99").Value = 2154.7
$ws.Range("M99").Value = -650.3332999999998
$ws.Range("N99").Value = -5150.7

$ws.Range("H126").Value = 2152.3125
$ws.Range("I126").Value = 2148.3333
$ws.Range("J126").Value = 2154.7
$ws.Range("K126").Value = 6444.999899999999
$ws.Range("L126").Value = 6464.099999999999
$ws.Range("M126").Value = -3974.999899999999
$ws.Range("N126").Value = -11404.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 738.7091
$ws.Range("J131").Value = 1306.2632
$ws.Range("L131").Value = 3918.7896
$ws.Range("N131").Value = -13998.7896

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 468212.1
$ws.Range("I102").Value = 1242.25
$ws.Range("J102").Value = 1869121.6
$ws.Range("K102").Value = 1242.25
$ws.Range("L102").Value = 1869121.6
$ws.Range("M102").Value = 379.75
$ws.Range("N102").Value = -1872365.6

$ws.Range("H126").Value = 465668.53
$ws.Range("I126").Value = 1669.5
$ws.Range("K126").Value = 5008.5
$ws.Range("M126").Value = -2538.5

$ws.Range("H132").Value = 1900.0646
$ws.Range("I132").Value = 1514.8889
$ws.Range("K132").Value = 4544.6667
$ws.Range("M132").Value = -2014.6667

$ws.Range("H141").Value = 22809.666
$ws.Range("J141").Value = 22809.666
$ws.Range("L141").Value = 22809.666
$ws.Range("N141").Value = -33169.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2659.7
$ws.Range("I7").Value = 1839.4
$ws.Range("J7").Value = 3480
$ws.Range("K7").Value = 1839.4
$ws.Range("L7").Value = 3480
$ws.Range("M7").Value = -1727.4
$ws.Range("N7").Value = -3704

$ws.Range("H20").Value = 5969.125
$ws.Range("J20").Value = 5969.125
$ws.Range("L20").Value = 5969.125
$ws.Range("N20").Value = -6421.125

$ws.Range("H46").Value = 2303.889
$ws.Range("I46").Value = 850
$ws.Range("J46").Value = 2863.077
$ws.Range("K46").Value = 850
$ws.Range("L46").Value = 2863.077
$ws.Range("M46").Value = -662
$ws.Range("N46").Value = -3239.077

$ws.Range("H68").Value = 1606.6666
$ws.Range("I68").Value = 992.3077
$ws.Range("J68").Value = 5600
$ws.Range("K68").Value = 992.3077
$ws.Range("L68").Value = 5600
$ws.Range("M68").Value = -243.3077
$ws.Range("N68").Value = -7098

$ws.Range("H71").Value = 1606.6666
$ws.Range("I71").Value = 992.3077
$ws.Range("J71").Value = 5600
$ws.Range("K71").Value = 4961.5385
$ws.Range("L71").Value = 28000
$ws.Range("M71").Value = -1217.5385
$ws.Range("N71").Value = -35488

$ws.Range("H100").Value = 2532.75
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 2532.75
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 2532.75
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -3614.75

$ws.Range("H122").Value = 2668.3333
$ws.Range("I122").Value = 2355.6099
$ws.Range("J122").Value = 3654.6155
$ws.Range("K122").Value = 7066.8297
$ws.Range("L122").Value = 10963.8465
$ws.Range("M122").Value = -4616.8297
$ws.Range("N122").Value = -15863.8465

$ws.Range("H126").Value = 2659.7
$ws.Range("I126").Value = 1839.4
$ws.Range("J126").Value = 3480
$ws.Range("K126").Value = 5518.200000000001
$ws.Range("L126").Value = 10440
$ws.Range("M126").Value = -3048.200000000001
$ws.Range("N126").Value = -15380

$ws.Range("H132").Value = 2220.2632
$ws.Range("I132").Value = 1715.4166
$ws.Range("J132").Value = 3085.7144
$ws.Range("K132").Value = 5146.2498
$ws.Range("L132").Value = 9257.143199999999
$ws.Range("M132").Value = -2616.2498
$ws.Range("N132").Value = -14317.1432

$ws.Range("H133").Value = 29493.334
$ws.Range("J133").Value = 29493.334
$ws.Range("L133").Value = 29493.334
$ws.Range("N133").Value = -34553.334

$ws.Range("H135").Value = 29470
$ws.Range("J135").Value = 29470
$ws.Range("L135").Value = 29470
$ws.Range("N135").Value = -39610

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 614.2143
$ws.Range("I81").Value = 574.9167
$ws.Range("K81").Value = 1149.8334
$ws.Range("M81").Value = -88.83339999999998

$ws.Range("H84").Value = 614.2143
$ws.Range("I84").Value = 574.9167
$ws.Range("K84").Value = 5749.166999999999
$ws.Range("M84").Value = -445.1669999999995

$ws.Range("H96").Value = 7434.3125
$ws.Range("I96").Value = 1291.6666
$ws.Range("J96").Value = 11119.9
$ws.Range("K96").Value = 1291.6666
$ws.Range("L96").Value = 11119.9
$ws.Range("M96").Value = 81.33339999999998
$ws.Range("N96").Value = -13865.9

$ws.Range("H140").Value = 40000
$ws.Range("J140").Value = 40000
$ws.Range("L140").Value = 40000
$ws.Range("N140").Value = -50360
